$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1a"
$ws.Range("C2").Value = "Il1rap"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09179766666666667
$ws.Range("H2").Value = 0.275393
$ws.Range("I2").Value = 0.01443540132615123
$ws.Range("J2").Value = 0.01443540132615123
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.652098
$ws.Range("N2").Value = 7.956294000000001
$ws.Range("O2").Value = 0.1783170651171435
$ws.Range("P2").Value = 0.1783170651171436
$ws.Range("Q2").Value = 0.2434564081713333
$ws.Range("R2").Value = 2.191107673542
$ws.Range("S2").Value = 0.00257407839826741
$ws.Range("T2").Value = 0.00257407839826741

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1a"
$ws.Range("C3").Value = "Il1rap"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09179766666666667
$ws.Range("H3").Value = 0.275393
$ws.Range("I3").Value = 0.01443540132615123
$ws.Range("J3").Value = 0.01443540132615123
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.511188000000001
$ws.Range("N3").Value = 16.533564
$ws.Range("O3").Value = 0.370551491486672
$ws.Range("P3").Value = 0.3705514914866722
$ws.Range("Q3").Value = 0.5059141989613334
$ws.Range("R3").Value = 4.553227790652
$ws.Range("S3").Value = 0.005349059491614024
$ws.Range("T3").Value = 0.005349059491614024

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1a"
$ws.Range("C4").Value = "Il1rap"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09179766666666667
$ws.Range("H4").Value = 0.275393
$ws.Range("I4").Value = 0.01443540132615123
$ws.Range("J4").Value = 0.01443540132615123
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.28691
$ws.Range("N4").Value = 6.860729999999999
$ws.Range("O4").Value = 0.1537632015811809
$ws.Range("P4").Value = 0.1537632015811809
$ws.Range("Q4").Value = 0.2099330018766666
$ws.Range("R4").Value = 1.88939701689
$ws.Range("S4").Value = 0.002219633524018238
$ws.Range("T4").Value = 0.002219633524018238

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il1a"
$ws.Range("C5").Value = "Il1rap"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09179766666666667
$ws.Range("H5").Value = 0.275393
$ws.Range("I5").Value = 0.01443540132615123
$ws.Range("J5").Value = 0.01443540132615123
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.422738333333334
$ws.Range("N5").Value = 13.268215
$ws.Range("O5").Value = 0.2973682418150034
$ws.Range("P5").Value = 0.2973682418150034
$ws.Range("Q5").Value = 0.4059970592772223
$ws.Range("R5").Value = 3.653973533495
$ws.Range("S5").Value = 0.004292629912251561
$ws.Range("T5").Value = 0.004292629912251561

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Il1a"
$ws.Range("C6").Value = "Il1rap"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.267406666666666
$ws.Range("H6").Value = 18.80222
$ws.Range("I6").Value = 0.9855645986738488
$ws.Range("J6").Value = 0.9855645986738487
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.652098
$ws.Range("N6").Value = 7.956294000000001
$ws.Range("O6").Value = 0.1783170651171435
$ws.Range("P6").Value = 0.1783170651171436
$ws.Range("Q6").Value = 16.62177668585333
$ws.Range("R6").Value = 149.59599017268
$ws.Range("S6").Value = 0.1757429867188761
$ws.Range("T6").Value = 0.1757429867188761

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Il1a"
$ws.Range("C7").Value = "Il1rap"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.267406666666666
$ws.Range("H7").Value = 18.80222
$ws.Range("I7").Value = 0.9855645986738488
$ws.Range("J7").Value = 0.9855645986738487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.511188000000001
$ws.Range("N7").Value = 16.533564
$ws.Range("O7").Value = 0.370551491486672
$ws.Range("P7").Value = 0.3705514914866722
$ws.Range("Q7").Value = 34.54085641245334
$ws.Range("R7").Value = 310.86770771208
$ws.Range("S7").Value = 0.365202431995058
$ws.Range("T7").Value = 0.3652024319950581

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Il1a"
$ws.Range("C8").Value = "Il1rap"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.267406666666666
$ws.Range("H8").Value = 18.80222
$ws.Range("I8").Value = 0.9855645986738488
$ws.Range("J8").Value = 0.9855645986738487
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.28691
$ws.Range("N8").Value = 6.860729999999999
$ws.Range("O8").Value = 0.1537632015811809
$ws.Range("P8").Value = 0.1537632015811809
$ws.Range("Q8").Value = 14.33299498006666
$ws.Range("R8").Value = 128.9969548206
$ws.Range("S8").Value = 0.1515435680571627
$ws.Range("T8").Value = 0.1515435680571627

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Il1a"
$ws.Range("C9").Value = "Il1rap"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.267406666666666
$ws.Range("H9").Value = 18.80222
$ws.Range("I9").Value = 0.9855645986738488
$ws.Range("J9").Value = 0.9855645986738487
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.422738333333334
$ws.Range("N9").Value = 13.268215
$ws.Range("O9").Value = 0.2973682418150034
$ws.Range("P9").Value = 0.2973682418150034
$ws.Range("Q9").Value = 27.71909971525556
$ws.Range("R9").Value = 249.4718974373
$ws.Range("S9").Value = 0.2930756119027518
$ws.Range("T9").Value = 0.2930756119027518
